$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.00226073850791258
$ws.Range("C2").Value = 0.00301431801055011
$ws.Range("D2").Value = 0.014318010550113
$ws.Range("E2").Value = 0.957045968349661
$ws.Range("F2").Value = 0.0188394875659382
$ws.Range("G2").Value = 0.975131876412962
$ws.Range("H2").Value = 0.0158251695553881
$ws.Range("I2").Value = 0.012810851544838
$ws.Range("J2").Value = 0.969856819894499
$ws.Range("K2").Value = 0.00678221552373775
$ws.Range("L2").Value = 0.972117558402411
$ws.Range("M2").Value = 0.00226073850791258
$ws.Range("N2").Value = 0.0105501130369254
$ws.Range("O2").Value = 0.0504898266767144
$ws.Range("P2").Value = 0.00452147701582517
$ws.Range("Q2").Value = 1
$ws.Range("R2").Value = 0.00602863602110023
$ws.Range("S2").Value = 0.000753579502637528
$ws.Range("T2").Value = 0.000753579502637528
$ws.Range("U2").Value = 0.00452147701582517
$ws.Range("V2").Value = 0.0113036925395629
$ws.Range("W2").Value = 0.00452147701582517
$ws.Range("X2").Value = 0.0173323285606631
$ws.Range("B3").Value = 0.0188394875659382
$ws.Range("C3").Value = 0.0278824415975885
$ws.Range("D3").Value = 0.0052750565184627
$ws.Range("E3").Value = 0.00226073850791258
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 0.00376789751318764
$ws.Range("H3").Value = 0.978146194423512
$ws.Range("I3").Value = 0.970610399397136
$ws.Range("J3").Value = 0.0233609645817634
$ws.Range("K3").Value = 0.027128862094951
$ws.Range("L3").Value = 0.025621703089676
$ws.Range("M3").Value = 0.00226073850791258
$ws.Range("N3").Value = 0.859080633006782
$ws.Range("O3").Value = 0.0497362471740769
$ws.Range("P3").Value = 0
$ws.Range("Q3").Value = 0
$ws.Range("R3").Value = 0.0195930670685757
$ws.Range("S3").Value = 0.99698568198945
$ws.Range("T3").Value = 0.995478522984175
$ws.Range("U3").Value = 0.0497362471740769
$ws.Range("V3").Value = 0
$ws.Range("W3").Value = 0.0120572720422005
$ws.Range("X3").Value = 0
$ws.Range("B4").Value = 0.970610399397136
$ws.Range("C4").Value = 0.00753579502637528
$ws.Range("D4").Value = 0.0120572720422005
$ws.Range("E4").Value = 0.0293896006028636
$ws.Range("F4").Value = 0.966842501883949
$ws.Range("G4").Value = 0.0195930670685757
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0.00301431801055011
$ws.Range("J4").Value = 0.00226073850791258
$ws.Range("K4").Value = 0.961567445365486
$ws.Range("L4").Value = 0.00226073850791258
$ws.Range("M4").Value = 0.000753579502637528
$ws.Range("N4").Value = 0.00301431801055011
$ws.Range("O4").Value = 0.00376789751318764
$ws.Range("P4").Value = 0.995478522984175
$ws.Range("Q4").Value = 0
$ws.Range("R4").Value = 0.972117558402411
$ws.Range("S4").Value = 0
$ws.Range("T4").Value = 0.000753579502637528
$ws.Range("U4").Value = 0.000753579502637528
$ws.Range("V4").Value = 0.974378296910324
$ws.Range("W4").Value = 0.970610399397136
$ws.Range("X4").Value = 0.968349660889224
$ws.Range("B5").Value = 0.00753579502637528
$ws.Range("C5").Value = 0.961567445365486
$ws.Range("D5").Value = 0.968349660889224
$ws.Range("E5").Value = 0.0113036925395629
$ws.Range("F5").Value = 0.014318010550113
$ws.Range("G5").Value = 0.00150715900527506
$ws.Range("H5").Value = 0.00602863602110023
$ws.Range("I5").Value = 0.0135644310474755
$ws.Range("J5").Value = 0.00452147701582517
$ws.Range("K5").Value = 0.00452147701582517
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = 0.994724943481537
$ws.Range("N5").Value = 0.127354935945742
$ws.Range("O5").Value = 0.896006028636021
$ws.Range("P5").Value = 0
$ws.Range("Q5").Value = 0
$ws.Range("R5").Value = 0.00150715900527506
$ws.Range("S5").Value = 0.00226073850791258
$ws.Range("T5").Value = 0.00301431801055011
$ws.Range("U5").Value = 0.94498869630746
$ws.Range("V5").Value = 0.0135644310474755
$ws.Range("W5").Value = 0.0120572720422005
$ws.Range("X5").Value = 0.014318010550113
